$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 value
$ws.Range("A2").Value = 917514074672

# Remove row 4 (A4 had value 919343865508) - delete the entire row
$ws.Range("A4").EntireRow.Delete()

# Update the active selection to B8
$ws.Range("B8").Select()
